$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster table (rows 2-19), as it should appear after the edit:
#  - Row 2 "Chris Paul" replaced by "Naji Marshall" / Dallas Mavericks
#  - Row order reshuffled
#  - Old "Corey Kispert" row replaced by "Jrue Holiday" / Boston Celtics
#  - New row 19 appended: "Jerami Grant" / Portland Trail Blazers
$data = @(
    @("Naji Marshall", "SG,SF", "Dallas Mavericks"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jrue Holiday", "PG,SG", "Boston Celtics"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
